$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Wipe existing values but keep the existing cell formatting
#    (bold/bordered/centred style on header row + column A, etc.)
# ------------------------------------------------------------------
$ws.Range("A1:M8").ClearContents()

# We need 12 data rows (rows 2-13) in the end, but only rows 2-8 (7 rows)
# currently carry the correctly-styled column A / row formatting.
# Clone the formatting of row 2 down onto the 5 extra rows we need (9-13).
$ws.Range("A2:M2").Copy() | Out-Null
$ws.Range("A9:M13").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 2. Header row (A1:L1) - unchanged apart from the removed "subject" column
# ------------------------------------------------------------------
$header = New-Object 'object[,]' 1,12
$header[0,0]  = "option"
$header[0,1]  = "max mark"
$header[0,2]  = "'9"
$header[0,3]  = "'8"
$header[0,4]  = "'7"
$header[0,5]  = "'6"
$header[0,6]  = "'5"
$header[0,7]  = "'4"
$header[0,8]  = "'3"
$header[0,9]  = "'2"
$header[0,10] = "'1"
$header[0,11] = "date"
$ws.Range("A1:L1").Value = $header

# ------------------------------------------------------------------
# 3. Data rows (A2:L13)
#    Columns: A option, B max mark, C..K grade thresholds (9..1), L date
# ------------------------------------------------------------------
$data = @(
    @("AY 12,32", 160,   0,   0,   0,   0,  93,  79,  58, 38, "18", "June 2022"),
    @("BY 22,42", 200, 168, 147, 127, 105,  84,  63,  40,  0,    0, "June 2022"),
    @("AX 11,31", 160,   0,   0,   0,   0, 113,  82,  60, 38, "16", "November 2022"),
    @("BX 21,41", 200, 165, 146, 127, 109,  91,  74,  54,  0,    0, "November 2022"),
    @("AY 12,32", 160,   0,   0,   0,   0, 113,  94,  71, 49, "27", "June 2023"),
    @("BY 22,42", 200, 181, 164, 148, 124, 101,  78,  55,  0,    0, "June 2023"),
    @("AY 12,32", 160,   0,   0,   0,   0, 109,  89,  66, 43, "21", "June 2024"),
    @("BY 22,42", 200, 180, 165, 150, 128, 106,  85,  60,  0,    0, "June 2024"),
    @("AY 12,32", 160,   0,   0,   0,   0, 103,  86,  63, 40, "17", "June 2025"),
    @("BY 22,42", 200, 184, 168, 152, 130, 108,  86,  65,  0,    0, "June 2025"),
    @("AX 02,11", 200, 158, 147, 136, 126, 116, 106,  84, 62, "40", "November 2024"),
    @("BX 11,82", 200, 158, 147, 136, 126, 116, 106,  84, 62, "40", "November 2024")
)

$rowCount = $data.Count
$startRow = 2

for ($r = 0; $r -lt $rowCount; $r++) {
    $rowNum = $startRow + $r
    $rowVals = $data[$r]

    # Column A: option label
    $ws.Cells.Item($rowNum, 1).Value = $rowVals[0]

    # Columns B..J: numeric mark columns (max mark, then grades 9..2)
    for ($c = 1; $c -le 9; $c++) {
        $ws.Cells.Item($rowNum, $c + 1).Value = $rowVals[$c]
    }

    # Column K (grade 1 threshold) - numeric 0 in some rows, text in others
    $kVal = $rowVals[10]
    $kCell = $ws.Cells.Item($rowNum, 11)
    if ($kVal -is [string]) {
        $kCell.Value = "'" + $kVal
    } else {
        $kCell.Value = $kVal
    }

    # Column L ("date" label like "June 2022") - must stay text, not become a date
    $lCell = $ws.Cells.Item($rowNum, 12)
    $lCell.Value = "'" + $rowVals[11]
}

# ------------------------------------------------------------------
# 4. Remove the now unused "subject" column entirely (M)
# ------------------------------------------------------------------
$ws.Columns("M:M").Delete() | Out-Null

# ------------------------------------------------------------------
# 5. Re-assert the original cell formatting (the text-forcing steps above
#    can attach stray number-format / quote-prefix styles to cells) so the
#    final look matches the original template exactly: bold+bordered+
#    centred style on the header row and column A, default style elsewhere.
# ------------------------------------------------------------------
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A1:L1").PasteSpecial(-4122) | Out-Null
$ws.Range("A2:A13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("B2").Copy() | Out-Null
$ws.Range("B2:L13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
